# Quarterly data refresh: the "Trimestre" (C) / "Valor" (D) series for each
# of the three region blocks (Brasil, Nordeste, Sergipe) rolls forward by one
# quarter - each row now shows what used to be in the next row, and a new
# quarter (01/01/2024) is appended at the end of every block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-QuarterCell($ws, $row, $text) {
    # Force the cell to stay text (dd/mm/yyyy-looking strings otherwise get
    # auto-parsed into date serials by COM `.Value` assignment), then drop
    # the temporary number-format back to the sheet's default style so no
    # stray formatting is left behind.
    $cell = $ws.Cells.Item($row, 3)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 52, 53, 54, 55, 56, 57, 58, 59, 60, 61, 62, 63, 64)

$newTrimestre = @("01/01/2019", "01/04/2019", "01/07/2019", "01/10/2019", "01/01/2020", "01/04/2020", "01/07/2020", "01/10/2020", "01/01/2021", "01/04/2021", "01/07/2021", "01/10/2021", "01/01/2022", "01/04/2022", "01/07/2022", "01/10/2022", "01/01/2023", "01/04/2023", "01/07/2023", "01/10/2023", "01/01/2024", "01/01/2019", "01/04/2019", "01/07/2019", "01/10/2019", "01/01/2020", "01/04/2020", "01/07/2020", "01/10/2020", "01/01/2021", "01/04/2021", "01/07/2021", "01/10/2021", "01/01/2022", "01/04/2022", "01/07/2022", "01/10/2022", "01/01/2023", "01/04/2023", "01/07/2023", "01/10/2023", "01/01/2024", "01/01/2019", "01/04/2019", "01/07/2019", "01/10/2019", "01/01/2020", "01/04/2020", "01/07/2020", "01/10/2020", "01/01/2021", "01/04/2021", "01/07/2021", "01/10/2021", "01/01/2022", "01/04/2022", "01/07/2022", "01/10/2022", "01/01/2023", "01/04/2023", "01/07/2023", "01/10/2023", "01/01/2024")

$newValor = @(87.15383963941923, 87.85947559951479, 88.09875854372994, 88.91898936863468, 87.6269256467444, 86.40199837580568, 85.10970347929863, 85.82012456093744, 85.09170501959174, 85.7680202656022, 87.35882739828995, 88.85372779747212, 88.85603969260613, 90.6967300113522, 91.29946932281176, 92.05869818976858, 91.20616836197172, 91.96054185222719, 92.31096399578379, 92.59072488218143, 92.07634205061291, 84.55102285920174, 85.24212449847059, 85.39254559873116, 86.24122476500217, 84.21958837190678, $null, $null, $null, $null, $null, $null, $null, $null, 87.30259251929546, 88.02800283174703, 89.13957176843775, 87.76160329045526, 88.67195362505535, 89.1498039836851, 89.56449309852451, 88.8911025222138, 84.47789275634995, 84.6503178928247, 85.21897810218978, 85.06666666666666, 84.21052631578947, $null, $null, $null, $null, $null, $null, $null, $null, 87.24954462659382, 87.87037037037037, 88.04744525547446, 88.1740775780511, 89.76303317535546, 90.20332717190388, 88.70214752567693, 89.92805755395683)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]

    Set-QuarterCell $ws $r $newTrimestre[$i]

    $dCell = $ws.Cells.Item($r, 4)
    if ($newValor[$i] -eq $null) {
        $dCell.ClearContents()
    } else {
        $dCell.Value = $newValor[$i]
    }
}
